$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 7) of data after the existing 6 data rows.
$row = 7

$ws.Cells.Item($row, 1).Value = (Get-Date -Year 2016 -Month 8 -Day 30 -Hour 21 -Minute 29 -Second 10)
$ws.Cells.Item($row, 2).Value = -2
$ws.Cells.Item($row, 3).Value = 50
$ws.Cells.Item($row, 4).Value = 46
$ws.Cells.Item($row, 5).Value = 20
$ws.Cells.Item($row, 6).Value = 80
$ws.Cells.Item($row, 7).Value = 9640
$ws.Cells.Item($row, 8).Value = 16674
$ws.Cells.Item($row, 9).Value = 1786
$ws.Cells.Item($row, 10).Value = 267
$ws.Cells.Item($row, 11).Value = 243
$ws.Cells.Item($row, 12).Value = 1
$ws.Cells.Item($row, 13).Value = 4
$ws.Cells.Item($row, 14).Value = "Bag"
